# [Npc] CompassUI를 위한 데이터 추가2
# Insert a new "Priority" (int) column into the MapNpcMenu sheet, right
# before the existing "bool" / ShowBubble column, shifting the remaining
# columns one to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MapNpcMenu")

# Insert a new column at H; everything from H onward (bool/ShowBubble,
# FunctionType/enum, Order/int) shifts right by one column.
$ws.Columns.Item(8).Insert()

# Populate the newly inserted column with the Priority field metadata.
$ws.Range("H1").Value = "int"
$ws.Range("H2").Value = "Priority"
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 10

# Match the author's final selection.
$ws.Range("I7").Select() | Out-Null
